$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 25000
$ws.Range("J3").Value = 25000
$ws.Range("L3").Value = 25000
$ws.Range("N3").Value = -25228
$ws.Range("H11").Value = 78815.92999999999
$ws.Range("I11").Value = 78815.92999999999
$ws.Range("K11").Value = 78815.92999999999
$ws.Range("M11").Value = -78675.92999999999
$ws.Range("H12").Value = 254.7
$ws.Range("I12").Value = 135.42857
$ws.Range("K12").Value = 135.42857
$ws.Range("M12").Value = 34.57142999999999
$ws.Range("H18").Value = 1882
$ws.Range("I18").Value = 1882
$ws.Range("K18").Value = 1882
$ws.Range("M18").Value = -1598
$ws.Range("H76").Value = 107019.9
$ws.Range("I76").Value = 7480
$ws.Range("K76").Value = 7480
$ws.Range("M76").Value = -7165
$ws.Range("H79").Value = 107019.9
$ws.Range("I79").Value = 7480
$ws.Range("K79").Value = 7480
$ws.Range("M79").Value = -6388
$ws.Range("H100").Value = 4659.0713
$ws.Range("I100").Value = 5318.9165
$ws.Range("K100").Value = 5318.9165
$ws.Range("M100").Value = -4777.9165
$ws.Range("H102").Value = 25000
$ws.Range("J102").Value = 25000
$ws.Range("L102").Value = 25000
$ws.Range("N102").Value = -31490
$ws.Range("H107").Value = 999.5
$ws.Range("I107").Value = 999.5
$ws.Range("K107").Value = 999.5
$ws.Range("M107").Value = 920.5
$ws.Range("H121").Value = 1990.3462
$ws.Range("J121").Value = 1990.3462
$ws.Range("L121").Value = 5971.0386
$ws.Range("N121").Value = -9465.0386
$ws.Range("H135").Value = 43479484
$ws.Range("I135").Value = 45455732
$ws.Range("K135").Value = 409101588
$ws.Range("M135").Value = -409099053

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2455.3845
$ws.Range("I32").Value = 2493.64
$ws.Range("K32").Value = 2493.64
$ws.Range("M32").Value = -2206.64
$ws.Range("H61").Value = 41669536
$ws.Range("I61").Value = 47621612
$ws.Range("K61").Value = 47621612
$ws.Range("M61").Value = -47621400
$ws.Range("H122").Value = 10419688
$ws.Range("I122").Value = 13336041
$ws.Range("K122").Value = 40008123
$ws.Range("M122").Value = -40005673
$ws.Range("H132").Value = 27779708
$ws.Range("I132").Value = 29413662
$ws.Range("K132").Value = 88240986
$ws.Range("M132").Value = -88238456
$ws.Range("H136").Value = 41669536
$ws.Range("I136").Value = 47621612
$ws.Range("K136").Value = 142864836
$ws.Range("M136").Value = -142862286

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1752.25
$ws.Range("I134").Value = 1462.3636
$ws.Range("K134").Value = 4387.0908
$ws.Range("M134").Value = -1852.0908

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 696
$ws.Range("I22").Value = 745
$ws.Range("K22").Value = 745
$ws.Range("M22").Value = -395
$ws.Range("H31").Value = 2806.8386
$ws.Range("I31").Value = 1523.8636
$ws.Range("K31").Value = 1523.8636
$ws.Range("M31").Value = -1228.8636
$ws.Range("H34").Value = 2806.8386
$ws.Range("I34").Value = 1523.8636
$ws.Range("K34").Value = 1523.8636
$ws.Range("M34").Value = -1321.8636
$ws.Range("H62").Value = 166673840
$ws.Range("J62").Value = 333338660
$ws.Range("L62").Value = 333338660
$ws.Range("N62").Value = -333339908
$ws.Range("H65").Value = 166673840
$ws.Range("J65").Value = 333338660
$ws.Range("L65").Value = 1666693300
$ws.Range("N65").Value = -1666699540
$ws.Range("H105").Value = 1567.3334
$ws.Range("I105").Value = 1480.8
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1480.8
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 266.2
$ws.Range("N105").Value = -5494
$ws.Range("H134").Value = 1902.2941
$ws.Range("I134").Value = 1502.9231
$ws.Range("J134").Value = 3200.25
$ws.Range("K134").Value = 4508.7693
$ws.Range("L134").Value = 9600.75
$ws.Range("M134").Value = -1973.7693
$ws.Range("N134").Value = -14670.75
$ws.Range("H141").Value = 85897.60000000001
$ws.Range("J141").Value = 85897.60000000001
$ws.Range("L141").Value = 85897.60000000001
$ws.Range("N141").Value = -96257.60000000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 631
$ws.Range("I51").Value = 631
$ws.Range("K51").Value = 1893
$ws.Range("M51").Value = -1433
$ws.Range("H87").Value = 342.33334
$ws.Range("I87").Value = 342.33334
$ws.Range("K87").Value = 1027.00002
$ws.Range("M87").Value = 220.9999800000001
$ws.Range("H90").Value = 342.33334
$ws.Range("I90").Value = 342.33334
$ws.Range("K90").Value = 3081.00006
$ws.Range("M90").Value = 3158.99994
$ws.Range("H122").Value = 855.2222
$ws.Range("J122").Value = 930.8333
$ws.Range("L122").Value = 8377.4997
$ws.Range("N122").Value = -13277.4997
$ws.Range("H131").Value = 4258.4414
$ws.Range("I131").Value = 852.5454999999999
$ws.Range("J131").Value = 5887.3477
$ws.Range("K131").Value = 2557.6365
$ws.Range("L131").Value = 17662.0431
$ws.Range("M131").Value = 2482.3635
$ws.Range("N131").Value = -27742.0431
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -32060

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1991.1
$ws.Range("I102").Value = 1587.2858
$ws.Range("K102").Value = 1587.2858
$ws.Range("M102").Value = 34.71419999999989
$ws.Range("H132").Value = 3924.394
$ws.Range("I132").Value = 3737.1904
$ws.Range("K132").Value = 11211.5712
$ws.Range("M132").Value = -8681.5712

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 351000
$ws.Range("J4").Value = 26500
$ws.Range("L4").Value = 26500
$ws.Range("N4").Value = -26726
$ws.Range("H7").Value = 83335336
$ws.Range("I7").Value = 125001500
$ws.Range("K7").Value = 125001500
$ws.Range("M7").Value = -125001388
$ws.Range("H22").Value = 768.2857
$ws.Range("I22").Value = 480
$ws.Range("J22").Value = 984.5
$ws.Range("K22").Value = 480
$ws.Range("L22").Value = 984.5
$ws.Range("M22").Value = -185
$ws.Range("N22").Value = -1574.5
$ws.Range("H27").Value = 768.2857
$ws.Range("I27").Value = 480
$ws.Range("J27").Value = 984.5
$ws.Range("K27").Value = 480
$ws.Range("L27").Value = 984.5
$ws.Range("M27").Value = -373
$ws.Range("N27").Value = -1198.5
$ws.Range("H28").Value = 351000
$ws.Range("J28").Value = 26500
$ws.Range("L28").Value = 26500
$ws.Range("N28").Value = -26964
$ws.Range("H37").Value = 351000
$ws.Range("J37").Value = 26500
$ws.Range("L37").Value = 26500
$ws.Range("N37").Value = -26714
$ws.Range("H46").Value = 2493.5625
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2493.5625
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2493.5625
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2869.5625
$ws.Range("H57").Value = 27666.334
$ws.Range("I57").Value = 18999
$ws.Range("J57").Value = 32000
$ws.Range("K57").Value = 18999
$ws.Range("L57").Value = 32000
$ws.Range("M57").Value = -18433
$ws.Range("N57").Value = -33132
$ws.Range("H126").Value = 83335336
$ws.Range("I126").Value = 125001500
$ws.Range("K126").Value = 375004500
$ws.Range("M126").Value = -375002030
$ws.Range("H132").Value = 5755.6333
$ws.Range("I132").Value = 3511.4707
$ws.Range("J132").Value = 8690.308000000001
$ws.Range("K132").Value = 10534.4121
$ws.Range("L132").Value = 26070.924
$ws.Range("M132").Value = -8004.4121
$ws.Range("N132").Value = -31130.924
$ws.Range("H136").Value = 2047
$ws.Range("I136").Value = 1868.5
$ws.Range("K136").Value = 5605.5
$ws.Range("M136").Value = -3055.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 30000
$ws.Range("J27").Value = 30000
$ws.Range("L27").Value = 30000
$ws.Range("N27").Value = -30138
$ws.Range("H54").Value = 7500
$ws.Range("H95").Value = 39171.75
$ws.Range("J95").Value = 39171.75
$ws.Range("L95").Value = 39171.75
$ws.Range("N95").Value = -44663.75
$ws.Range("H97").Value = 10000
$ws.Range("J97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("N97").Value = -11982
$ws.Range("H109").Value = 88800
$ws.Range("J109").Value = 88800
$ws.Range("L109").Value = 88800
$ws.Range("N109").Value = -91574
$ws.Range("H122").Value = 2494.625
$ws.Range("I122").Value = 2492.8333
$ws.Range("K122").Value = 7478.499899999999
$ws.Range("M122").Value = -5028.499899999999
$ws.Range("H132").Value = 2604.3
$ws.Range("I132").Value = 1867.0646
$ws.Range("K132").Value = 5601.1938
$ws.Range("M132").Value = -3071.1938
